$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.116.71"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "1.638.94"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.994"
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.77"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.994"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.257"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0634"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.79"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "1.863.46"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "1.633.29"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.32"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "26.061.12"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.994"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.33"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.99"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.37"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.994"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.87"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.89"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0496"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.24"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.909"
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("D37").Value = "1.144.86"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.547"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.50"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.994"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.26"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.797"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "1.772.85"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.73"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  -4.80%  "
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("E49").Value = "  +4.65%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.62"
$ws.Range("E51").Value = "  +1.08%  "
